# 4/12/18 mods to file
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Row 20 ("Category 38"): these counts no longer roll up into E/F/G/H,
# so clear the old "=C20" link (was only in G20) and zero out the row's
# E/F/G/H cells outright.
$ws.Range("E20").Value = 0
$ws.Range("F20").Value = 0
$ws.Range("G20").Value = 0
$ws.Range("H20").Value = 0

# --- Row 21 ("Category 39"): fill in the missing F/G/H links to C21
# (E21 already had "=C21").
$ws.Range("F21").Formula = "=C21"
$ws.Range("G21").Formula = "=C21"
$ws.Range("H21").Formula = "=C21"

# --- Row 27 ("Category 45"): fill in the missing E/F links to C27
# (G27/H27 already had the formula).
$ws.Range("E27").Formula = "=C27"
$ws.Range("F27").Formula = "=C27"

# --- Row 29 ("Category 47"): same treatment as row 20 -- drop the
# "=C29" links (previously in E/F/G/H) and zero the row out.
$ws.Range("E29").Value = 0
$ws.Range("F29").Value = 0
$ws.Range("G29").Value = 0
$ws.Range("H29").Value = 0

# --- Sheet view: selection moved from C22 to I29, and the window no
# longer needs to be scrolled to keep row 7 pinned at the top.
$ws.Range("I29").Select()
